# Generate Report for Handback
# Adds a new handback record (a5368211-9e58-402a-b5ee-7f35c9e9b558) as row 4
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$guid = "a5368211-9e58-402a-b5ee-7f35c9e9b558"
$mdName = "$guid.md"
$inSync = "Handed back: in sync with en-US"
$include = "Include"

function Add-Link($ws, $addr, $url, $disp) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $disp) | Out-Null
}

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

$overviewMdUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handback/OpenLocalizationTest/oltest/xinjiang/" + $mdName
Add-Link $wsOverview "A4" $overviewMdUrl $mdName

# ---------------------------------------------------------------------------
# zh-cn sheet:
# Source File Name | Status | Correspond Handoff File | Correspond Handoff Datetime |
# Target File | Correspond Handback File | Correspond Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfName = $guid + ".d37b0ab66537e39e9ada6663e41a75b661b59990.zh-cn.xlf"
$zhHandoffDate = "2016-01-28 03:56:39"
$zhHandbackDate = "2016-01-28 03:57:20"

$wsZhCn.Range("B4").Value = $inSync
$wsZhCn.Range("D4").Value = $zhHandoffDate
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("G4").Value = $zhHandbackDate
$wsZhCn.Range("H4").Value = $include

$zhMdUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handback/OpenLocalizationTest/oltest/xinjiang/" + $mdName
$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/" + $zhXlfName
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/" + $zhXlfName

Add-Link $wsZhCn "A4" $zhMdUrl $mdName
Add-Link $wsZhCn "C4" $zhHandoffUrl $zhXlfName
Add-Link $wsZhCn "E4" $zhMdUrl $mdName
Add-Link $wsZhCn "F4" $zhHandbackUrl $zhXlfName

# ---------------------------------------------------------------------------
# de-de sheet: same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfName = $guid + ".d37b0ab66537e39e9ada6663e41a75b661b59990.de-de.xlf"
$deHandoffDate = "2016-01-28 03:56:50"
$deHandbackDate = "2016-01-28 03:57:39"

$wsDeDe.Range("B4").Value = $inSync
$wsDeDe.Range("D4").Value = $deHandoffDate
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("G4").Value = $deHandbackDate
$wsDeDe.Range("H4").Value = $include

$deMdUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handback/OpenLocalizationTest/oltest/xinjiang/" + $mdName
$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/" + $deXlfName
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/" + $deXlfName

Add-Link $wsDeDe "A4" $deMdUrl $mdName
Add-Link $wsDeDe "C4" $deHandoffUrl $deXlfName
Add-Link $wsDeDe "E4" $deMdUrl $mdName
Add-Link $wsDeDe "F4" $deHandbackUrl $deXlfName

Write-Output "Handback row added for $guid"
